$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Row 13 gains a "siftelc" / reservoir filter on the Storage process set and
# switches its AndOr combinator from AND to Or (per chat with Antti).
$ws.Range("A13").Value = "ELE,STG,NST"
$ws.Range("B13").Value = "*siftelc*"
$ws.Range("C13").Value = "*pondage*,*large reservoir*"
$ws.Range("H13").Value = "Or"
$ws.Range("I13").Value = "Or"
